$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    "B2" = 0.9828542969916043; "D2" = 0.2790060750674499; "E2" = 0.2628485844196806; "F2" = 1.302266585538554; "G2" = 0.652325958620132; "H2" = 0.7727522898643855; "J2" = 0.3174931185742711; "K2" = 0.3335988923132902; "L2" = 0.1276347857931341; "N2" = 1.916304565947325; "O2" = 2.8414267891999
    "B3" = 0.9580547559902755; "D3" = 0.278785218416651; "E3" = 0.264423178977264; "F3" = 1.306937131406613; "G3" = 0.6542885991542775; "H3" = 0.7770238932122311; "J3" = 0.3198960946071354; "K3" = 0.2909906278841277; "L3" = 0.1185356960727404; "N3" = 1.933150029894559; "O3" = 2.85428989063206
    "B4" = 0.9432028571903288; "D4" = 0.2787456415265339; "E4" = 0.2654715403786536; "F4" = 1.310390589078438; "G4" = 0.6558534163982159; "H4" = 0.7799290108576145; "J4" = 0.3214592418078386; "K4" = 0.2647389308204993; "L4" = 0.1129727470878095; "N4" = 1.944061902427041; "O4" = 2.863533033390681
    "B5" = 0.9372456096342887; "D5" = 0.278753743811798; "E5" = 0.2659193022151278; "F5" = 1.311945382481511; "G5" = 0.6565816047809818; "H5" = 0.7811839649513033; "J5" = 0.3221183247255919; "K5" = 0.2540192306992992; "L5" = 0.1107119981586067; "N5" = 1.94865167845947; "O5" = 2.86763824758971
    "B6" = 0.9362621727317446; "D6" = 0.2787565550637296; "E6" = 0.2659948947939128; "F6" = 1.312212468024335; "G6" = 0.6567079881492717; "H6" = 0.7813966460143362; "J6" = 0.3222290999376112; "K6" = 0.2522379300091586; "L6" = 0.1103369821801721; "N6" = 1.949422453034426; "O6" = 2.868340372509095
    "B7" = 0.9431221300049231; "D7" = 0.2787456525789835; "E7" = 0.2654774958095185; "F7" = 1.310410960158542; "G7" = 0.6558628704513012; "H7" = 0.7799456476141557; "J7" = 0.3214680409548301; "K7" = 0.264594448992284; "L7" = 0.1129422325022915; "N7" = 1.944123222216433; "O7" = 2.863587026522538
    "B8" = 0.9742259665208621; "D8" = 0.2789100416739814; "E8" = 0.2633746012924778; "F8" = 1.303755539026618; "G8" = 0.6529280408670033; "H8" = 0.7741665992227382; "J8" = 0.3183034800431601; "K8" = 0.3189266741529764; "L8" = 0.1244925588156462; "N8" = 1.921994910343169; "O8" = 2.845582984602984
    "B9" = 1.038170166422731; "D9" = 0.2799909705117045; "E9" = 0.2598962250277719; "F9" = 1.295343891537136; "G9" = 0.6500256469634991; "H9" = 0.7650698913564753; "J9" = 0.3127922229783078; "K9" = 0.4247305855384411; "L9" = 0.147325768263201; "N9" = 1.88310876670689; "O9" = 2.820938202118924
    "B10" = 1.086917821664144; "D10" = 0.2812433604755284; "E10" = 0.2577318198588809; "F10" = 1.291982908361405; "G10" = 0.6496310870635966; "H10" = 0.7597442012310864; "J10" = 0.3091643025903323; "K10" = 0.5019837514771268; "L10" = 0.16420543944038; "N10" = 1.857279505018294; "O10" = 2.809317132823622
    "B11" = 1.109472157096661; "D11" = 0.2819117356934413; "E11" = 0.2568316289921828; "F11" = 1.291064132245609; "G11" = 0.6498287231699891; "H11" = 0.7576150773379737; "J11" = 0.307604869192506; "K11" = 0.5370180480360887; "L11" = 0.1719055547347637; "N11" = 1.84612258162532; "O11" = 2.80543590072574
    "B12" = 1.118066730227866; "D12" = 0.2821789396567169; "E12" = 0.2565028494929127; "F12" = 1.290803777618962; "G12" = 0.649957757837015; "H12" = 0.7568509533788443; "J12" = 0.3070273942970125; "K12" = 0.5502683652751159; "L12" = 0.1748243114716246; "N12" = 1.841982916813667; "O12" = 2.804167988624386
    "B13" = 1.11621335688352; "D13" = 0.282120766493847; "E13" = 0.2565731202941777; "F13" = 1.290855957890031; "G13" = 0.6499275581897166; "H13" = 0.7570136487861703; "J13" = 0.3071511839217909; "K13" = 0.5474154149846413; "L13" = 0.1741955795796031; "N13" = 1.842870678188358; "O13" = 2.804432083850543
    "B14" = 1.110178165053583; "D14" = 0.2819334365131425; "E14" = 0.2568043377351383; "F14" = 1.291040958755268; "G14" = 0.6498382530049156; "H14" = 0.7575513685537345; "J14" = 0.3075570986206948; "K14" = 0.5381084930853319; "L14" = 0.1721456258125897; "N14" = 1.845780301216404; "O14" = 2.805327545354061
    "B15" = 1.10648841202962; "D15" = 0.2818205259174391; "E15" = 0.2569475402713213; "F15" = 1.291165675650944; "G15" = 0.6497906076080824; "H15" = 0.7578862213826767; "J15" = 0.3078074316403541; "K15" = 0.5324055738034588; "L15" = 0.1708903404813924; "N15" = 1.847573626233554; "O15" = 2.805902317767107
    "B16" = 1.08545140225155; "D16" = 0.2812016588035675; "E16" = 0.2577923449995883; "F16" = 1.292055215567572; "G16" = 0.6496257567291792; "H16" = 0.7598892444262049; "J16" = 0.3092680426341712; "K16" = 0.4996919249666121; "L16" = 0.1637026339619041; "N16" = 1.858020563818052; "O16" = 2.809599036607693
    "B17" = 1.072642372244957; "D17" = 0.2808472139184204; "E17" = 0.2583321996384331; "F17" = 1.292757070994426; "G17" = 0.649621191214564; "H17" = 0.7611931616103078; "J17" = 0.3101873496463234; "K17" = 0.4795947877655919; "L17" = 0.1592985672754423; "N17" = 1.86458125517931; "O17" = 2.81222659681049
    "B18" = 1.065310653599511; "D18" = 0.280652637532981; "E18" = 0.2586506566766626; "F18" = 1.293218200447804; "G18" = 0.6496540640695372; "H18" = 0.7619707765023946; "J18" = 0.3107246706520961; "K18" = 0.4680252760202563; "L18" = 0.156767497068401; "N18" = 1.868410603837887; "O18" = 2.813870183614142
    "B19" = 1.062834412891959; "D19" = 0.2805883555371338; "E19" = 0.2587598467034748; "F19" = 1.293384201735101; "G19" = 0.6496712922868539; "H19" = 0.7622388126676753; "J19" = 0.310908069293216; "K19" = 0.4641063158172187; "L19" = 0.1559108752626202; "N19" = 1.869716742725181; "O19" = 2.814449401226796
    "B20" = 1.074002226194722; "D20" = 0.2808839841983968; "E20" = 0.2582739089638686; "F20" = 1.292676413634588; "G20" = 0.6496180034830559; "H20" = 0.761051497767383; "J20" = 0.3100886021165365; "K20" = 0.4817352206805197; "L20" = 0.1597671787658186; "N20" = 1.863877082323164; "O20" = 2.811933199463681
    "B21" = 1.111949394656136; "D21" = 0.2819880777231987; "E21" = 0.2567360953638484; "F21" = 1.290984244384539; "G21" = 0.6498630135768337; "H21" = 0.7573922845244851; "J21" = 0.3074375177051749; "K21" = 0.5408426134054594; "L21" = 0.172747670056566; "N21" = 1.844923361893652; "O21" = 2.805059051072163
    "B22" = 1.137062935956521; "D22" = 0.2827918420702531; "E22" = 0.2558015763138446; "F22" = 1.290388634615425; "G22" = 0.6503390097941946; "H22" = 0.7552462984989177; "J22" = 0.3057809264513605; "K22" = 0.5793766869228989; "L22" = 0.1812479262641347; "N22" = 1.833032718488532; "O22" = 2.80174268019843
    "B23" = 1.123630980603082; "D23" = 0.2823553659247011; "E23" = 0.256293904498154; "F23" = 1.290659885723883; "G23" = 0.6500560717083061; "H23" = 0.7563692142319525; "J23" = 0.3066581308442959; "K23" = 0.5588193773593275; "L23" = 0.1767097115810685; "N23" = 1.839333543563079; "O23" = 2.803405141732725
    "B24" = 1.073387335360593; "D24" = 0.2808673316965127; "E24" = 0.2583002369842937; "F24" = 1.292712699322607; "G24" = 0.6496193340817484; "H24" = 0.7611154568425604; "J24" = 0.3101332184644221; "K24" = 0.4807675788677557; "L24" = 0.1595553167352932; "N24" = 1.864195259698352; "O24" = 2.812065430202978
    "B25" = 1.020558919923559; "D25" = 0.2796177765033647; "E25" = 0.2607683614505767; "F25" = 1.297123731011261; "G25" = 0.6505055434513736; "H25" = 0.7672919792986406; "J25" = 0.3142090425246273; "K25" = 0.3961901940039354; "L25" = 0.1411299560608654; "N25" = 1.893146730403689; "O25" = 2.826465353883492
}

foreach ($key in $data.Keys) {
    $ws.Range($key).Value = $data[$key]
}
